$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$s = $wb.Styles.Add("MyBoldStyle")
$s.Font.Bold = $true
$s.Font.Name = "Arial Unicode MS"
$s.Font.Size = 10
$s.HorizontalAlignment = -4108
$s.VerticalAlignment = -4108

$r = $ws.Range("BN33")
$r.Value = 10
$r.Style = "MyBoldStyle"
